# Update dashboards - 2026-02-13
# Refresh the rolling-window economic indicator figures on the
# "Aguilar Prototype" sheet: new release dates come in, the Q/R/S/T/U
# trailing-observations window shifts, and the "most recent release"
# highlight (yellow fill, style 49) moves to whichever date column is now
# the newest for that block of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Style moves: the yellow "latest" highlight shifts off of N18:N21
#    (CPI / Core CPI date cells) onto N51 (30y Mortgage date), while the
#    highlight on C22:C25 (Vehicle Sales / Consumer Credit date cells)
#    is removed (reverts to the plain/no-fill style).
#    Reuse PasteSpecial(formats) from cells that already carry the
#    target style so the workbook keeps reusing the existing cellXf
#    instead of fabricating a near-duplicate one.
# ---------------------------------------------------------------------

# Cells that should pick up the highlighted "latest date" style (49)
# (PasteSpecial only honors the first area of a multi-area Range in this
# host, so paste into each target cell individually.)
$ws.Range("N31").Copy()
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("N21").PasteSpecial(-4122)
$ws.Range("N51").PasteSpecial(-4122)

# Cells that should drop back to the plain "older date" style (48)
$ws.Range("N22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C25").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Release-date cells (new data pulled in)
# ---------------------------------------------------------------------
$ws.Range("N18").Value = 46023
$ws.Range("N19").Value = 46023
$ws.Range("N20").Value = 46023
$ws.Range("N21").Value = 46023

$ws.Range("N29").Value = 46065
$ws.Range("N30").Value = 46065

$ws.Range("C42").Value = 46023
$ws.Range("C43").Value = 46023

$ws.Range("N47").Value = 46064
$ws.Range("N48").Value = 46064
$ws.Range("N49").Value = 46064
$ws.Range("N50").Value = 46064

$ws.Range("N51").Value = 46062
$ws.Range("N52").Value = 46064

# ---------------------------------------------------------------------
# 3) CPI / Core CPI trailing windows (rows 18-21)
#    (S18/T19/S20/T21 are cleared explicitly - along with the other
#    already-blank cells in this block below - so they serialize back
#    out as clean empty cells rather than picking up a stray 0.)
# ---------------------------------------------------------------------
$ws.Range("Q18").Value = 0.00170842649932057
$ws.Range("R18").Value = 0.00297788428704604
$ws.Range("S18").ClearContents()
$ws.Range("T18").ClearContents()
$ws.Range("U18").Value = 0.002950901819104068

$ws.Range("Q19").Value = 0.02391201432150015
$ws.Range("R19").Value = 0.02653304114557758
$ws.Range("S19").Value = 0.02696443916493949
$ws.Range("T19").ClearContents()
$ws.Range("U19").Value = 0.03022571584713336

$ws.Range("Q20").Value = 0.002950448142634121
$ws.Range("R20").Value = 0.002329002576704653
$ws.Range("S20").ClearContents()
$ws.Range("T20").ClearContents()
$ws.Range("U20").Value = 0.002177737336973129

$ws.Range("Q21").Value = 0.02512028782828883
$ws.Range("R21").Value = 0.02646484707309002
$ws.Range("S21").Value = 0.02599044806094405
$ws.Range("T21").ClearContents()
$ws.Range("U21").Value = 0.03019966825885779

# ---------------------------------------------------------------------
# 4) 5yr5yr Forward / 10yr TIPS breakeven (rows 29-30)
# ---------------------------------------------------------------------
$ws.Range("Q29").Value = 2.13
$ws.Range("R29").Value = 2.15
$ws.Range("S29").Value = 2.17
$ws.Range("T29").Value = 2.2
$ws.Range("U29").Value = 2.18

$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.32
$ws.Range("S30").Value = 2.32
$ws.Range("T30").Value = 2.35
$ws.Range("U30").Value = 2.34

# ---------------------------------------------------------------------
# 5) Existing Home Sales level + Y/Y delta (rows 42-43)
# ---------------------------------------------------------------------
$ws.Range("F42").Value = 3910000
$ws.Range("G42").Value = 4270000
$ws.Range("H42").Value = 4090000
$ws.Range("I42").Value = 4110000
$ws.Range("J42").Value = 4080000
$ws.Range("S42").ClearContents()

$ws.Range("F43").Value = -0.04400977995110025
$ws.Range("G43").ClearContents()
$ws.Range("H43").ClearContents()
$ws.Range("I43").ClearContents()
$ws.Range("J43").ClearContents()
$ws.Range("R43").ClearContents()
$ws.Range("S43").ClearContents()

# ---------------------------------------------------------------------
# 6) Rates block (rows 48-52): 2y/5y/10y UST, 30y Mortgage, BAA
# ---------------------------------------------------------------------
$ws.Range("Q48").Value = 3.52
$ws.Range("R48").Value = 3.45
$ws.Range("S48").Value = 3.48
$ws.Range("T48").Value = 3.5
$ws.Range("U48").Value = 3.47

$ws.Range("Q49").Value = 3.75
$ws.Range("R49").Value = 3.7
$ws.Range("S49").Value = 3.75
$ws.Range("T49").Value = 3.76
$ws.Range("U49").Value = 3.74

$ws.Range("Q50").Value = 4.18
$ws.Range("R50").Value = 4.16
$ws.Range("T50").Value = 4.22
$ws.Range("U50").Value = 4.21

$ws.Range("Q51").Value = 6.09
$ws.Range("R51").Value = 6.11
$ws.Range("S51").Value = 6.1
$ws.Range("T51").Value = 6.09
$ws.Range("U51").Value = 6.06

$ws.Range("Q52").Value = 5.85
$ws.Range("R52").Value = 5.82
$ws.Range("S52").Value = 5.86
$ws.Range("T52").Value = 5.87
$ws.Range("U52").Value = 5.88
